# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 gets a new built-in table style (tableStyleId).
# 2) The design/theme applied to the deck is swapped from "Integral" to
#    the default "Office Theme" colour palette (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), which is what happens when a new Design is picked
#    from the Design gallery.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{16B9FFDC-0221-443A-84CF-19C4F6013F85}")
    }
}

# --- 2) Swap the applied design's colour scheme to the Office Theme ------
# (RGB long values = R + G*256 + B*65536, same encoding VBA's RGB() returns)
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
